$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 527. Existing rows 527-565 shift down
# to become rows 529-567 (new price observations are prepended for this
# market/week, the rest of the weekly series shifts down).
$ws.Rows.Item(527).Insert()
$ws.Rows.Item(527).Insert()

# --- New row 527 ---
$ws.Cells.Item(527, 1).Value = 3
$ws.Cells.Item(527, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(527, 3).Value = "Coquimbo"
$ws.Cells.Item(527, 4).Value = 44610
$ws.Cells.Item(527, 5).Value = 5
$ws.Cells.Item(527, 6).Value = "Fruta"
$ws.Cells.Item(527, 7).Value = 100109
$ws.Cells.Item(527, 8).Value = "Uva"
$ws.Cells.Item(527, 9).Value = 100109001
$ws.Cells.Item(527, 10).Value = "Uva"
$ws.Cells.Item(527, 11).Value = "Red Globe"
$ws.Cells.Item(527, 12).Value = "Primera"
$ws.Cells.Item(527, 13).Value = 85
$ws.Cells.Item(527, 14).Value = 12000
$ws.Cells.Item(527, 15).Value = 12000
$ws.Cells.Item(527, 16).Value = 12000
$ws.Cells.Item(527, 17).Value = "`$/caja 12 kilos"
$ws.Cells.Item(527, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(527, 19).Value = 1000
$ws.Cells.Item(527, 20).Value = 12

# --- New row 528 ---
$ws.Cells.Item(528, 1).Value = 3
$ws.Cells.Item(528, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(528, 3).Value = "Coquimbo"
$ws.Cells.Item(528, 4).Value = 44610
$ws.Cells.Item(528, 5).Value = 5
$ws.Cells.Item(528, 6).Value = "Fruta"
$ws.Cells.Item(528, 7).Value = 100109
$ws.Cells.Item(528, 8).Value = "Uva"
$ws.Cells.Item(528, 9).Value = 100109001
$ws.Cells.Item(528, 10).Value = "Uva"
$ws.Cells.Item(528, 11).Value = "Superior Seedless"
$ws.Cells.Item(528, 12).Value = "Primera"
$ws.Cells.Item(528, 13).Value = 75
$ws.Cells.Item(528, 14).Value = 15000
$ws.Cells.Item(528, 15).Value = 15000
$ws.Cells.Item(528, 16).Value = 15000
$ws.Cells.Item(528, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(528, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(528, 19).Value = 1000
$ws.Cells.Item(528, 20).Value = 15
